$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 6
    3 = 4
    4 = 5
    5 = 6
    6 = 7
    7 = 3
    8 = 3
    9 = 3
    10 = 4
    11 = 2
    12 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
